# Regenerate quadratic/linear problem data (new random sample)
# Mirrors the author's "volver a generar problemas cuadraticos y lineales" commit:
# values on several sheets are replaced with a freshly generated sample.

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "-0.25 - x + y_1 + y_2"
$ws.Range("B2").Value = "0.25"
$ws.Range("D2").Value = "0.34"
$ws.Range("E2").Value = "4.8"
$ws.Range("F2").Value = "3.5999999999999996"

$ws.Range("A3").Value = "1.7000000000000002 - y_1"
$ws.Range("B3").Value = "-1.7000000000000002"
$ws.Range("D3").Value = "0.14"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "2.3000000000000003"

$ws.Range("A4").Value = "-3.8 - y_2"
$ws.Range("B4").Value = "-3.8"
$ws.Range("D4").Value = "0.38"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0.4"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "5.25"
$ws.Range("B2").Value = "1.7000000000000002"
$ws.Range("C2").Value = "3.8"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A2").Value = "-1.9000000000000001"
$ws.Range("A3").Value = "-0.96"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item("Vector_BF")
$ws.Range("A2").Value = "4.8"
$ws.Range("A3").Value = "-4.5"
$ws.Range("A4").Value = "-6.6"
